# Update the cryptocurrency listing with refreshed prices / volume figures.
# (Generated to match the "Updated cryptos list ... with GitHub Actions" commit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds locale-formatted text (e.g. "68.130.86", "0.0000289")
# rather than real numbers, so force it to stay plain text before writing the
# refreshed values - otherwise COM auto-coerces numeric-looking text to a
# binary double and we lose the original formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.130.86"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "3.620.79"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "587.48"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").Value = "194.07"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("D7").Value = "3.616.50"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "0.680"
$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("D12").Value = "55.70"
$ws.Range("E12").Value = "  -2.83%  "

$ws.Range("D13").Value = "0.0000289"
$ws.Range("E13").Value = "  +6.32%  "

$ws.Range("D14").Value = "10.02"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").Value = "4.195.72"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").Value = "3.620.98"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "68.129.74"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "18.54"
$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("E21").Value = "  -2.44%  "

$ws.Range("D22").Value = "406.03"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "13.46"
$ws.Range("E23").Value = "  +23.03%  "

$ws.Range("D24").Value = "4.27"
$ws.Range("E24").Value = "  -3.71%  "

$ws.Range("D25").Value = "86.08"
$ws.Range("E25").Value = "  -2.79%  "

$ws.Range("E26").Value = "  +0.31%  "

# Row 27: coin re-ranked
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "3.95"
$ws.Range("E27").Value = "  +6.96%  "

# Row 28: coin re-ranked
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "12.60"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").Value = "6.13"
$ws.Range("E29").Value = "  +1.00%  "

$ws.Range("D30").Value = "8.19"
$ws.Range("E30").Value = "  +14.32%  "

$ws.Range("D31").Value = "9.16"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").Value = "31.59"
$ws.Range("E32").Value = "  -1.14%  "

$ws.Range("D33").Value = "677.93"
$ws.Range("E33").Value = "  +11.70%  "

$ws.Range("D34").Value = "12.26"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("E35").Value = "  +1.04%  "

$ws.Range("D36").Value = "64.56"
$ws.Range("E36").Value = "  -5.27%  "

$ws.Range("D37").Value = "42.63"
$ws.Range("E37").Value = "  -3.63%  "

$ws.Range("D38").Value = "0.423"
$ws.Range("E38").Value = "  +7.74%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("D40").Value = "0.0₃0788"
$ws.Range("E40").Value = "  +2.02%  "

$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +17.25%  "

# Row 42: coin re-ranked
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.215.33"
$ws.Range("E42").Value = "  +15.80%  "

# Row 43: coin re-ranked
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "3.13"
$ws.Range("E43").Value = "  +7.82%  "

$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").Value = "0.0423"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("E47").Value = "  -2.66%  "

# Row 48: coin re-ranked
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "8.82"
$ws.Range("E48").Value = "  -1.21%  "

# Row 49: coin re-ranked
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.11"
$ws.Range("E49").Value = "  -3.50%  "

$ws.Range("D50").Value = "143.40"
$ws.Range("E50").Value = "  -0.21%  "

# Row 51: coin re-ranked
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E51").Value = "  -1.19%  "

